# Update LR-pair NATMI statistics for Thbs2-Cd47 with recomputed TPM-based values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3.0
$ws.Range("F2").Value = 1.0
$ws.Range("G2").Value = 1.281134
$ws.Range("H2").Value = 3.843402
$ws.Range("I2").Value = 0.007312702338676299
$ws.Range("J2").Value = 0.007312702338676299
$ws.Range("M2").Value = 57.35848733333334
$ws.Range("N2").Value = 172.075462
$ws.Range("O2").Value = 0.261658309594631
$ws.Range("P2").Value = 0.261658309594631
$ws.Range("Q2").Value = 73.48390831130267
$ws.Range("R2").Value = 661.355174801724
$ws.Range("S2").Value = 0.001913429332506745
$ws.Range("T2").Value = 0.001913429332506745

# Row 3
$ws.Range("E3").Value = 3.0
$ws.Range("F3").Value = 1.0
$ws.Range("G3").Value = 1.281134
$ws.Range("H3").Value = 3.843402
$ws.Range("I3").Value = 0.007312702338676299
$ws.Range("J3").Value = 0.007312702338676299
$ws.Range("O3").Value = 0.2957894889638607
$ws.Range("P3").Value = 0.2957894889638607
$ws.Range("Q3").Value = 83.06928115579868
$ws.Range("R3").Value = 747.6235304021881
$ws.Range("S3").Value = 0.002163020487701891
$ws.Range("T3").Value = 0.002163020487701891

# Row 4
$ws.Range("E4").Value = 3.0
$ws.Range("F4").Value = 1.0
$ws.Range("G4").Value = 1.281134
$ws.Range("H4").Value = 3.843402
$ws.Range("I4").Value = 0.007312702338676299
$ws.Range("J4").Value = 0.007312702338676299
$ws.Range("M4").Value = 29.294891
$ws.Range("N4").Value = 87.88467299999999
$ws.Range("O4").Value = 0.1336376186888105
$ws.Range("P4").Value = 0.1336376186888105
$ws.Range("Q4").Value = 37.53068088639399
$ws.Range("R4").Value = 337.776127977546
$ws.Range("S4").Value = 0.000977252126720796
$ws.Range("T4").Value = 0.000977252126720796

# Row 5
$ws.Range("E5").Value = 3.0
$ws.Range("F5").Value = 1.0
$ws.Range("G5").Value = 1.281134
$ws.Range("H5").Value = 3.843402
$ws.Range("I5").Value = 0.007312702338676299
$ws.Range("J5").Value = 0.007312702338676299
$ws.Range("M5").Value = 67.71760166666667
$ws.Range("N5").Value = 203.152805
$ws.Range("O5").Value = 0.3089145827526977
$ws.Range("P5").Value = 0.3089145827526977
$ws.Range("Q5").Value = 86.75532189362333
$ws.Range("R5").Value = 780.7978970426101
$ws.Range("S5").Value = 0.002259000391746866
$ws.Range("T5").Value = 0.002259000391746866

# Row 6
$ws.Range("I6").Value = 0.9398544320918915
$ws.Range("J6").Value = 0.9398544320918915
$ws.Range("M6").Value = 57.35848733333334
$ws.Range("N6").Value = 172.075462
$ws.Range("O6").Value = 0.261658309594631
$ws.Range("P6").Value = 0.261658309594631
$ws.Range("Q6").Value = 9444.412436772802
$ws.Range("R6").Value = 84999.71193095524
$ws.Range("S6").Value = 0.2459207219661862
$ws.Range("T6").Value = 0.2459207219661862

# Row 7
$ws.Range("I7").Value = 0.9398544320918915
$ws.Range("J7").Value = 0.9398544320918915
$ws.Range("O7").Value = 0.2957894889638607
$ws.Range("P7").Value = 0.2957894889638607
$ws.Range("S7").Value = 0.2779990621688801
$ws.Range("T7").Value = 0.2779990621688801

# Row 8
$ws.Range("I8").Value = 0.9398544320918915
$ws.Range("J8").Value = 0.9398544320918915
$ws.Range("M8").Value = 29.294891
$ws.Range("N8").Value = 87.88467299999999
$ws.Range("O8").Value = 0.1336376186888105
$ws.Range("P8").Value = 0.1336376186888105
$ws.Range("Q8").Value = 4823.576174288643
$ws.Range("R8").Value = 43412.18556859779
$ws.Range("S8").Value = 0.1255999082188847
$ws.Range("T8").Value = 0.1255999082188847

# Row 9
$ws.Range("I9").Value = 0.9398544320918915
$ws.Range("J9").Value = 0.9398544320918915
$ws.Range("M9").Value = 67.71760166666667
$ws.Range("N9").Value = 203.152805
$ws.Range("O9").Value = 0.3089145827526977
$ws.Range("P9").Value = 0.3089145827526977
$ws.Range("Q9").Value = 11150.10156478487
$ws.Range("R9").Value = 100350.9140830639
$ws.Range("S9").Value = 0.2903347397379403
$ws.Range("T9").Value = 0.2903347397379403

# Row 10
$ws.Range("G10").Value = 9.213772333333333
$ws.Range("H10").Value = 27.641317
$ws.Range("I10").Value = 0.05259213672418158
$ws.Range("J10").Value = 0.05259213672418158
$ws.Range("M10").Value = 57.35848733333334
$ws.Range("N10").Value = 172.075462
$ws.Range("O10").Value = 0.261658309594631
$ws.Range("P10").Value = 0.261658309594631
$ws.Range("Q10").Value = 528.4880436737171
$ws.Range("R10").Value = 4756.392393063455
$ws.Range("S10").Value = 0.01376116959321907
$ws.Range("T10").Value = 0.01376116959321907

# Row 11
$ws.Range("G11").Value = 9.213772333333333
$ws.Range("H11").Value = 27.641317
$ws.Range("I11").Value = 0.05259213672418158
$ws.Range("J11").Value = 0.05259213672418158
$ws.Range("O11").Value = 0.2957894889638607
$ws.Range("P11").Value = 0.2957894889638607
$ws.Range("Q11").Value = 597.4249723004665
$ws.Range("R11").Value = 5376.824750704199
$ws.Range("S11").Value = 0.01555620124516316
$ws.Range("T11").Value = 0.01555620124516316

# Row 12
$ws.Range("G12").Value = 9.213772333333333
$ws.Range("H12").Value = 27.641317
$ws.Range("I12").Value = 0.05259213672418158
$ws.Range("J12").Value = 0.05259213672418158
$ws.Range("M12").Value = 29.294891
$ws.Range("N12").Value = 87.88467299999999
$ws.Range("O12").Value = 0.1336376186888105
$ws.Range("P12").Value = 0.1336376186888105
$ws.Range("Q12").Value = 269.9164562038156
$ws.Range("R12").Value = 2429.248105834341
$ws.Range("S12").Value = 0.007028287913575965
$ws.Range("T12").Value = 0.007028287913575965

# Row 13
$ws.Range("G13").Value = 9.213772333333333
$ws.Range("H13").Value = 27.641317
$ws.Range("I13").Value = 0.05259213672418158
$ws.Range("J13").Value = 0.05259213672418158
$ws.Range("M13").Value = 67.71760166666667
$ws.Range("N13").Value = 203.152805
$ws.Range("O13").Value = 0.3089145827526977
$ws.Range("P13").Value = 0.3089145827526977
$ws.Range("Q13").Value = 623.9345647160205
$ws.Range("R13").Value = 5615.411082444185
$ws.Range("S13").Value = 0.01624647797222338
$ws.Range("T13").Value = 0.01624647797222338

# Row 14
$ws.Range("E14").Value = 1.0
$ws.Range("F14").Value = 0.3333333333333333
$ws.Range("G14").Value = 0.042174
$ws.Range("H14").Value = 0.126522
$ws.Range("I14").Value = 0.0002407288452506406
$ws.Range("J14").Value = 0.0002407288452506406
$ws.Range("M14").Value = 57.35848733333334
$ws.Range("N14").Value = 172.075462
$ws.Range("O14").Value = 0.261658309594631
$ws.Range("P14").Value = 0.261658309594631
$ws.Range("Q14").Value = 2.419036844796
$ws.Range("R14").Value = 21.771331603164
$ws.Range("S14").Value = 0.00006298870271895014
$ws.Range("T14").Value = 0.00006298870271895014

# Row 15
$ws.Range("E15").Value = 1.0
$ws.Range("F15").Value = 0.3333333333333333
$ws.Range("G15").Value = 0.042174
$ws.Range("H15").Value = 0.126522
$ws.Range("I15").Value = 0.0002407288452506406
$ws.Range("J15").Value = 0.0002407288452506406
$ws.Range("O15").Value = 0.2957894889638607
$ws.Range("P15").Value = 0.2957894889638607
$ws.Range("Q15").Value = 2.734580351052
$ws.Range("R15").Value = 24.611223159468
$ws.Range("S15").Value = 0.00007120506211554729
$ws.Range("T15").Value = 0.00007120506211554729

# Row 16
$ws.Range("E16").Value = 1.0
$ws.Range("F16").Value = 0.3333333333333333
$ws.Range("G16").Value = 0.042174
$ws.Range("H16").Value = 0.126522
$ws.Range("I16").Value = 0.0002407288452506406
$ws.Range("J16").Value = 0.0002407288452506406
$ws.Range("M16").Value = 29.294891
$ws.Range("N16").Value = 87.88467299999999
$ws.Range("O16").Value = 0.1336376186888105
$ws.Range("P16").Value = 0.1336376186888105
$ws.Range("Q16").Value = 1.235482733034
$ws.Range("R16").Value = 11.119344597306
$ws.Range("S16").Value = 0.00003217042962900278
$ws.Range("T16").Value = 0.00003217042962900278

# Row 17
$ws.Range("E17").Value = 1.0
$ws.Range("F17").Value = 0.3333333333333333
$ws.Range("G17").Value = 0.042174
$ws.Range("H17").Value = 0.126522
$ws.Range("I17").Value = 0.0002407288452506406
$ws.Range("J17").Value = 0.0002407288452506406
$ws.Range("M17").Value = 67.71760166666667
$ws.Range("N17").Value = 203.152805
$ws.Range("O17").Value = 0.3089145827526977
$ws.Range("P17").Value = 0.3089145827526977
$ws.Range("Q17").Value = 2.85592213269
$ws.Range("R17").Value = 25.70329919421
$ws.Range("S17").Value = 0.00007436465078714039
$ws.Range("T17").Value = 0.00007436465078714039

Write-Host "Applied NATMI TPM update to rows 2-17"
